$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.014835047558478
$ws.Range("D2").Value = 0.04148214030968056
$ws.Range("E2").Value = 0.2927579589039091
$ws.Range("F2").Value = 1.236456993981065
$ws.Range("G2").Value = 1.11535828446145
$ws.Range("H2").Value = 1.06594749751585
$ws.Range("K2").Value = 0.3856721015945084
$ws.Range("L2").Value = 0.1439342446123106
$ws.Range("B3").Value = 0.9907636128689603
$ws.Range("D3").Value = 0.04133848113034233
$ws.Range("E3").Value = 0.2939407682653901
$ws.Range("F3").Value = 1.218209807132737
$ws.Range("G3").Value = 1.096942095090753
$ws.Range("H3").Value = 1.062126300758635
$ws.Range("K3").Value = 0.3357844291435299
$ws.Range("L3").Value = 0.1338909479568997
$ws.Range("B4").Value = 0.976518745054733
$ws.Range("D4").Value = 0.04124646374188323
$ws.Range("E4").Value = 0.294734379954761
$ws.Range("F4").Value = 1.207686301774544
$ws.Range("G4").Value = 1.086274783360608
$ws.Range("H4").Value = 1.06025553162587
$ws.Range("K4").Value = 0.3051030483956936
$ws.Range("L4").Value = 0.1278025210259131
$ws.Range("B5").Value = 0.9708486908503744
$ws.Range("D5").Value = 0.04120801483429659
$ws.Range("E5").Value = 0.2950747382443701
$ws.Range("F5").Value = 1.203568621906896
$ws.Range("G5").Value = 1.082088305328483
$ws.Range("H5").Value = 1.059612631928402
$ws.Range("K5").Value = 0.2925876863056374
$ws.Range("L5").Value = 0.1253411061171477
$ws.Range("B6").Value = 0.9699153356748411
$ws.Range("D6").Value = 0.04120157318995155
$ws.Range("E6").Value = 0.2951322790054025
$ws.Range("F6").Value = 1.202895186200607
$ws.Range("G6").Value = 1.081402826173147
$ws.Range("H6").Value = 1.059513090606885
$ws.Range("K6").Value = 0.2905087728923093
$ws.Range("L6").Value = 0.1249335781632084
$ws.Range("B7").Value = 0.9764417303488244
$ws.Range("D7").Value = 0.0412459490474788
$ws.Range("E7").Value = 0.294738901472968
$ws.Range("F7").Value = 1.207630078468611
$ws.Range("G7").Value = 1.08621767367012
$ws.Range("H7").Value = 1.060246377735822
$ws.Range("K7").Value = 0.3049343120373749
$ws.Range("L7").Value = 0.1277692458783548
$ws.Range("B8").Value = 1.006424323061765
$ws.Range("D8").Value = 0.04143340151506969
$ws.Range("E8").Value = 0.2931518204212828
$ws.Range("F8").Value = 1.230023944222125
$ws.Range("G8").Value = 1.108875196035143
$ws.Range("H8").Value = 1.064531200369203
$ws.Range("K8").Value = 0.368481319398569
$ws.Range("L8").Value = 0.1404550852265913
$ws.Range("B9").Value = 1.069457376358173
$ws.Range("D9").Value = 0.04177047422624547
$ws.Range("E9").Value = 0.2905733403527382
$ws.Range("F9").Value = 1.2793556116508
$ws.Range("G9").Value = 1.158411781266523
$ws.Range("H9").Value = 1.076712610597497
$ws.Range("K9").Value = 0.4926983044041435
$ws.Range("L9").Value = 0.1659537250190226
$ws.Range("B10").Value = 1.11834613595309
$ws.Range("D10").Value = 0.04199915144917199
$ws.Range("E10").Value = 0.2890034096473659
$ws.Range("F10").Value = 1.318933992454788
$ws.Range("G10").Value = 1.197957865798514
$ws.Range("H10").Value = 1.087977051123858
$ws.Range("K10").Value = 0.5837262815918223
$ws.Range("L10").Value = 0.1850705617980424
$ws.Range("B11").Value = 1.141146225536431
$ws.Range("D11").Value = 0.04209899126679417
$ws.Range("E11").Value = 0.2883594713061566
$ws.Range("F11").Value = 1.337670426804806
$ws.Range("G11").Value = 1.216641472739354
$ws.Range("H11").Value = 1.093606699410714
$ws.Range("K11").Value = 0.6250889436056184
$ws.Range("L11").Value = 0.1938514647523562
$ws.Range("B12").Value = 1.149860424272532
$ws.Range("D12").Value = 0.04213618998208801
$ws.Range("E12").Value = 0.288125712357548
$ws.Range("F12").Value = 1.344871158981931
$ws.Range("G12").Value = 1.223816847556918
$ws.Range("H12").Value = 1.095811337768282
$ws.Range("K12").Value = 0.6407452402773686
$ws.Range("L12").Value = 0.197188758829185
$ws.Range("B13").Value = 1.147980098695569
$ws.Range("D13").Value = 0.0421282057342296
$ws.Range("E13").Value = 0.2881756081294018
$ws.Range("F13").Value = 1.343315647461253
$ws.Range("G13").Value = 1.222267033179037
$ws.Range("H13").Value = 1.095333288894324
$ws.Range("K13").Value = 0.6373736837709885
$ws.Range("L13").Value = 0.1964694721974638
$ws.Range("B14").Value = 1.141861540165536
$ws.Range("D14").Value = 0.04210206385851123
$ws.Range("E14").Value = 0.2883400377398306
$ws.Range("F14").Value = 1.338260715758651
$ws.Range("G14").Value = 1.217229782055398
$ws.Range("H14").Value = 1.093786616272467
$ws.Range("K14").Value = 0.6263771347064733
$ws.Range("L14").Value = 0.1941257823177125
$ws.Range("B15").Value = 1.138124195594088
$ws.Range("D15").Value = 0.04208597177825268
$ws.Range("E15").Value = 0.28844206881368
$ws.Range("F15").Value = 1.335178194480434
$ws.Range("G15").Value = 1.214157397439919
$ws.Range("H15").Value = 1.092848721077502
$ws.Range("K15").Value = 0.619640534423894
$ws.Range("L15").Value = 0.1926917876935335
$ws.Range("B16").Value = 1.11686734829479
$ws.Range("D16").Value = 0.04199254188146284
$ws.Range("E16").Value = 0.28904690468182
$ws.Range("F16").Value = 1.317724285734343
$ws.Range("G16").Value = 1.196750855299541
$ws.Range("H16").Value = 1.087619321974472
$ws.Range("K16").Value = 0.5810221855179236
$ws.Range("L16").Value = 0.184498410414534
$ws.Range("B17").Value = 1.103970280853929
$ws.Range("D17").Value = 0.04193414873178369
$ws.Range("E17").Value = 0.2894359307664605
$ws.Range("F17").Value = 1.307204643364372
$ws.Range("G17").Value = 1.18625059672118
$ws.Range("H17").Value = 1.084540796930469
$ws.Range("K17").Value = 0.557319066329967
$ws.Range("L17").Value = 0.1794936939808025
$ws.Range("B18").Value = 1.096604987741188
$ws.Range("D18").Value = 0.04190016884266257
$ws.Range("E18").Value = 0.2896662995110244
$ws.Range("F18").Value = 1.301222902399303
$ws.Range("G18").Value = 1.180276407477947
$ws.Range("H18").Value = 1.08281767293056
$ws.Range("K18").Value = 0.5436812994742013
$ws.Range("L18").Value = 0.1766230725597211
$ws.Range("B19").Value = 1.094120295964785
$ws.Range("D19").Value = 0.0418885964188469
$ws.Range("E19").Value = 0.2897454342787675
$ws.Range("F19").Value = 1.299209405979099
$ws.Range("G19").Value = 1.178264846096283
$ws.Range("H19").Value = 1.082242417438721
$ws.Range("K19").Value = 0.5390630424741687
$ws.Range("L19").Value = 0.1756524965490485
$ws.Range("B20").Value = 1.105337737367222
$ws.Range("D20").Value = 0.04194040557552015
$ws.Range("E20").Value = 0.289393834130065
$ws.Range("F20").Value = 1.30831734572385
$ws.Range("G20").Value = 1.18736160700314
$ws.Range("H20").Value = 1.08486358767172
$ws.Range("K20").Value = 0.5598427562954384
$ws.Range("L20").Value = 0.1800256308755053
$ws.Range("B21").Value = 1.143656532009686
$ws.Range("D21").Value = 0.04210975892414481
$ws.Range("E21").Value = 0.2882914671324883
$ws.Range("F21").Value = 1.33974260131609
$ws.Range("G21").Value = 1.218706618627493
$ws.Range("H21").Value = 1.094238934301131
$ws.Range("K21").Value = 0.629607274843778
$ws.Range("L21").Value = 0.1948138509343096
$ws.Range("B22").Value = 1.169168000658743
$ws.Range("D22").Value = 0.04221689174766752
$ws.Range("E22").Value = 0.2876297927275342
$ws.Range("F22").Value = 1.360896817114863
$ws.Range("G22").Value = 1.239777250963868
$ws.Range("H22").Value = 1.100790737654961
$ws.Range("K22").Value = 0.6751624440546493
$ws.Range("L22").Value = 0.2045496656444357
$ws.Range("B23").Value = 1.155509325785147
$ws.Range("D23").Value = 0.04216003977851912
$ws.Range("E23").Value = 0.2879775660335433
$ws.Range("F23").Value = 1.349549922970866
$ws.Range("G23").Value = 1.228477772902693
$ws.Range("H23").Value = 1.097255032202582
$ws.Range("K23").Value = 0.6508525053379799
$ws.Range("L23").Value = 0.1993469952039248
$ws.Range("B24").Value = 1.104719356372186
$ws.Range("D24").Value = 0.04193757812876697
$ws.Range("E24").Value = 0.2894128451107072
$ws.Range("F24").Value = 1.307814086982646
$ws.Range("G24").Value = 1.186859124452354
$ws.Range("H24").Value = 1.084717508319386
$ws.Range("K24").Value = 0.5587018288495642
$ws.Range("L24").Value = 0.1797851214727046
$ws.Range("B25").Value = 1.05195214277984
$ws.Range("D25").Value = 0.04168259306822719
$ws.Range("E25").Value = 0.2912138386490053
$ws.Range("F25").Value = 1.265426830870965
$ws.Range("G25").Value = 1.14446021004153
$ws.Range("H25").Value = 1.073011520537534
$ws.Range("K25").Value = 0.4591361077524709
$ws.Range("L25").Value = 0.1589887002936621
